# case 1: refresh the simulated random-walk values in columns A:B
# and bump the column width (15.42578125 -> 16.28515625 chars) to fit them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ColumnWidth is in "characters"; Excel stores width = ColumnWidth + 5/6 in the
# saved XML (rounded to the nearest 1/6 char), so 15.5 is the COM-property input
# that lands on the closest achievable stored width to 16.28515625.
$ws.Range("A:A").ColumnWidth = 15.5
$ws.Range("B:B").ColumnWidth = 15.5

# Refresh cell values in columns A and B, rows 1-32
$ws.Cells.Item(1, 1).Value2 = -0.08404241345678543
$ws.Cells.Item(1, 2).Value2 = 0.08361927445412931
$ws.Cells.Item(2, 1).Value2 = -0.07849591158427938
$ws.Cells.Item(2, 2).Value2 = 0.07759461429248393
$ws.Cells.Item(3, 1).Value2 = -0.02788253985340816
$ws.Cells.Item(3, 2).Value2 = 0.027699769204394542
$ws.Cells.Item(4, 1).Value2 = -0.019699769291849023
$ws.Cells.Item(4, 2).Value2 = 0.019548385495772536
$ws.Cells.Item(5, 1).Value2 = -0.016548385534271404
$ws.Cells.Item(5, 2).Value2 = 0.016054038224358358
$ws.Cells.Item(6, 1).Value2 = -0.008899363961074158
$ws.Cells.Item(6, 2).Value2 = 0.008856317376633172
$ws.Cells.Item(7, 1).Value2 = 0.0011436825160355824
$ws.Cells.Item(7, 2).Value2 = -0.0011452201362973646
$ws.Cells.Item(8, 1).Value2 = 0.011145220029479042
$ws.Cells.Item(8, 2).Value2 = -0.011154361215045316
$ws.Cells.Item(9, 1).Value2 = 0.013154361190590436
$ws.Cells.Item(9, 2).Value2 = -0.013172502868687097
$ws.Cells.Item(10, 1).Value2 = 0.015172502847530467
$ws.Cells.Item(10, 2).Value2 = -0.015172490565385388
$ws.Cells.Item(11, 1).Value2 = 0.018172490534557717
$ws.Cells.Item(11, 2).Value2 = -0.01817879430814351
$ws.Cells.Item(12, 1).Value2 = 0.021678794273428004
$ws.Cells.Item(12, 2).Value2 = -0.02178621821438531
$ws.Cells.Item(13, 1).Value2 = 0.025286218186552745
$ws.Cells.Item(13, 2).Value2 = -0.025378176061526148
$ws.Cells.Item(14, 1).Value2 = -0.009075958153919927
$ws.Cells.Item(14, 2).Value2 = 0.009049428366784795
$ws.Cells.Item(15, 1).Value2 = -0.008049428367497669
$ws.Cells.Item(15, 2).Value2 = 0.008032328826093504
$ws.Cells.Item(16, 1).Value2 = -0.0060323288373926864
$ws.Cells.Item(16, 2).Value2 = 0.006002962781551435
$ws.Cells.Item(17, 1).Value2 = -0.004002962793659748
$ws.Cells.Item(17, 2).Value2 = 0.003999999967744472
$ws.Cells.Item(18, 1).Value2 = 0.0008766370696449144
$ws.Cells.Item(18, 2).Value2 = -0.001056076233510339
$ws.Cells.Item(19, 1).Value2 = 0.005056076192838432
$ws.Cells.Item(19, 2).Value2 = -0.006455523783158501
$ws.Cells.Item(20, 1).Value2 = -0.0021456602969500693
$ws.Cells.Item(20, 2).Value2 = 0.0020581235254368835
$ws.Cells.Item(21, 1).Value2 = 0.001941876434694123
$ws.Cells.Item(21, 2).Value2 = -0.0020624633656112223
$ws.Cells.Item(22, 1).Value2 = -0.045712074485814824
$ws.Cells.Item(22, 2).Value2 = 0.04549914367395402
$ws.Cells.Item(23, 1).Value2 = -0.04049914373409891
$ws.Cells.Item(23, 2).Value2 = 0.040098892733897884
$ws.Cells.Item(24, 1).Value2 = -0.020098892949373948
$ws.Cells.Item(24, 2).Value2 = 0.0199999997813034
$ws.Cells.Item(25, 1).Value2 = -0.005154674292136008
$ws.Cells.Item(25, 2).Value2 = 0.005129591177810866
$ws.Cells.Item(26, 1).Value2 = -0.0026295912116260922
$ws.Cells.Item(26, 2).Value2 = 0.002599371305221254
$ws.Cells.Item(27, 1).Value2 = -0.00009937133902226947
$ws.Cells.Item(27, 2).Value2 = -0.00006977977952438152
$ws.Cells.Item(28, 1).Value2 = 0.002069779751318279
$ws.Cells.Item(28, 2).Value2 = -0.002169080011731772
$ws.Cells.Item(29, 1).Value2 = 0.009169079934610913
$ws.Cells.Item(29, 2).Value2 = -0.009189514086293649
$ws.Cells.Item(30, 1).Value2 = 0.06918951348428815
$ws.Cells.Item(30, 2).Value2 = -0.06950886764495001
$ws.Cells.Item(31, 1).Value2 = 0.00002249417602584458
$ws.Cells.Item(31, 2).Value2 = -0.00008587702082962778
$ws.Cells.Item(32, 1).Value2 = 0.010085876928094706
$ws.Cells.Item(32, 2).Value2 = -0.010122491249346055
